$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume-change (E) values for rows with simple in-place edits.
# A leading apostrophe forces Excel to keep the value as literal text (matching the
# original inline-string cells) instead of auto-converting numeric-looking text to a number.
$ws.Range("D2").Value = "'66.011.24"
$ws.Range("E2").Value = "'  -0.98%  "
$ws.Range("D3").Value = "'3.514.32"
$ws.Range("E3").Value = "'  +0.03%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.27%  "
$ws.Range("D5").Value = "'574.43"
$ws.Range("E5").Value = "'  +3.25%  "
$ws.Range("D6").Value = "'178.84"
$ws.Range("E6").Value = "'  -6.46%  "
$ws.Range("E7").Value = "'  +4.33%  "
$ws.Range("E9").Value = "'  -0.65%  "
$ws.Range("E10").Value = "'  +4.84%  "
$ws.Range("D11").Value = "'54.88"
$ws.Range("E11").Value = "'  -2.54%  "
$ws.Range("E12").Value = "'  +1.30%  "
$ws.Range("D13").Value = "'9.22"
$ws.Range("E13").Value = "'  -2.50%  "
$ws.Range("D14").Value = "'4.071.31"
$ws.Range("E14").Value = "'  -0.72%  "
$ws.Range("D15").Value = "'3.514.34"
$ws.Range("E15").Value = "'  -0.46%  "
$ws.Range("E16").Value = "'  +0.14%  "
$ws.Range("D17").Value = "'18.36"
$ws.Range("E17").Value = "'  +0.19%  "
$ws.Range("D18").Value = "'12.14"
$ws.Range("E18").Value = "'  +2.09%  "
$ws.Range("D19").Value = "'65.918.35"
$ws.Range("E19").Value = "'  -1.38%  "
$ws.Range("E20").Value = "'  +0.91%  "
$ws.Range("D21").Value = "'415.27"
$ws.Range("E21").Value = "'  +1.80%  "
$ws.Range("D22").Value = "'4.17"
$ws.Range("E22").Value = "'  +5.08%  "
$ws.Range("E23").Value = "'  +3.78%  "
$ws.Range("D24").Value = "'85.64"
$ws.Range("E24").Value = "'  +0.18%  "
$ws.Range("D25").Value = "'12.83"
$ws.Range("E25").Value = "'  +7.26%  "
$ws.Range("D26").Value = "'10.94"
$ws.Range("E26").Value = "'  -2.25%  "
$ws.Range("E27").Value = "'  -3.02%  "
$ws.Range("D28").Value = "'9.01"
$ws.Range("E28").Value = "'  +1.54%  "
$ws.Range("D29").Value = "'30.36"
$ws.Range("E29").Value = "'  -0.30%  "
$ws.Range("D30").Value = "'622.77"
$ws.Range("E30").Value = "'  -6.71%  "
$ws.Range("D31").Value = "'6.42"
$ws.Range("E31").Value = "'  -4.59%  "
$ws.Range("D32").Value = "'11.66"
$ws.Range("E32").Value = "'  -1.17%  "
$ws.Range("D33").Value = "'0.110"
$ws.Range("E33").Value = "'  -1.08%  "
$ws.Range("D34").Value = "'59.66"
$ws.Range("E34").Value = "'  -0.13%  "
$ws.Range("D35").Value = "'0.153"
$ws.Range("E35").Value = "'  +11.58%  "
$ws.Range("D36").Value = "'0.0₃0806"
$ws.Range("E36").Value = "'  -1.28%  "
$ws.Range("E37").Value = "'  +0.15%  "
$ws.Range("D38").Value = "'37.34"
$ws.Range("E38").Value = "'  -3.95%  "
$ws.Range("D41").Value = "'3.33"
$ws.Range("E41").Value = "'  -1.00%  "
$ws.Range("D42").Value = "'0.998"
$ws.Range("E42").Value = "'  -0.36%  "
$ws.Range("D43").Value = "'2.91"
$ws.Range("E43").Value = "'  -4.56%  "
$ws.Range("E44").Value = "'  +0.23%  "
$ws.Range("E45").Value = "'  -3.83%  "
$ws.Range("E46").Value = "'  -5.82%  "
$ws.Range("D47").Value = "'2.72"
$ws.Range("E47").Value = "'  -0.32%  "
$ws.Range("E48").Value = "'  +1.82%  "
$ws.Range("D49").Value = "'138.79"
$ws.Range("E49").Value = "'  +0.31%  "
$ws.Range("E50").Value = "'  -6.61%  "
$ws.Range("E51").Value = "'  -7.71%  "

# Rows 39 and 40 swap: TheGraph and Maker exchange row positions with updated data
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "'3.275.10"
$ws.Range("E39").Value = "'  +8.75%  "

$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "'0.380"
$ws.Range("E40").Value = "'  -3.75%  "
